$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the dates in column A for rows 6-9 so they all read 2016-03-03 (serial 42432)
$ws.Range("A6").Value = 42432
$ws.Range("A7").Value = 42432
$ws.Range("A8").Value = 42432
$ws.Range("A9").Value = 42432

# Update the selected range / active cell to match the corrected rows
$ws.Range("A7:A9").Select()
